$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.771.13"
$ws.Range("E2").Value = "  -0.23%  "
$ws.Range("D3").Value = "2.735.18"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.75%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.61"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").Value = "3.218.31"
$ws.Range("E13").Value = "  -0.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").Value = "63.612.56"
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "2.741.30"
$ws.Range("E17").Value = "  -0.19%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.42%  "
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.35%  "
$ws.Range("E25").Value = "  -0.11%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("D28").Value = "0.0₃0904"
$ws.Range("E28").Value = "  +2.06%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.39%  "
$ws.Range("E31").Value = "  +10.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.99"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.47"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("B36").Value = "USDe"
$ws.Range("C36").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "343.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.26"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.40%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.15%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "134.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.622"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.59%  "
$ws.Range("E48").Value = "  -0.66%  "
$ws.Range("E49").Value = "  -1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
